$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings remain stored as text (matching source formatting)
$textCells = @("D5", "D6", "D7", "D8", "D9", "D11", "D12", "D13", "D14", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply updated values from the crypto price refresh
$ws.Range("D2").Value = "24.901.54"
$ws.Range("E2").Value = "  +1.82%  "
$ws.Range("D3").Value = "1.672.46"
$ws.Range("E3").Value = "  +1.06%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "330.89"
$ws.Range("E5").Value = "  +7.62%  "
$ws.Range("D6").Value = "0.9976"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").Value = "0.3642"
$ws.Range("E7").Value = "  +0.47%  "
$ws.Range("D8").Value = "46.84"
$ws.Range("E8").Value = "  -1.37%  "
$ws.Range("D9").Value = "0.3235"
$ws.Range("E9").Value = "  -0.78%  "
$ws.Range("E10").Value = "  +1.54%  "
$ws.Range("D11").Value = "0.07057"
$ws.Range("E11").Value = "  +1.40%  "
$ws.Range("D12").Value = "0.9969"
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("D13").Value = "6.077"
$ws.Range("E13").Value = "  +2.66%  "
$ws.Range("D14").Value = "19.62"
$ws.Range("E14").Value = "  +1.67%  "
$ws.Range("D15").Value = "1.668.87"
$ws.Range("E15").Value = "  +0.80%  "
$ws.Range("D16").Value = "6.621"
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("D17").Value = "0.00001045"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").Value = "0.06546"
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("D19").Value = "0.9994"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").Value = "78.72"
$ws.Range("E20").Value = "  +3.10%  "
$ws.Range("D21").Value = "15.86"
$ws.Range("E21").Value = "  +0.97%  "
$ws.Range("D22").Value = "5.911"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "12.84"
$ws.Range("E23").Value = "  +2.56%  "
$ws.Range("D24").Value = "24.939.61"
$ws.Range("E24").Value = "  +2.02%  "
$ws.Range("D25").Value = "2.448"
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").Value = "2.398"
$ws.Range("E26").Value = "  +4.26%  "
$ws.Range("D27").Value = "148.23"
$ws.Range("E27").Value = "  +0.98%  "
$ws.Range("D28").Value = "18.72"
$ws.Range("E28").Value = "  +1.44%  "
$ws.Range("D29").Value = "1.849.57"
$ws.Range("E29").Value = "  +0.50%  "
$ws.Range("D30").Value = "125.69"
$ws.Range("E30").Value = "  +1.06%  "
$ws.Range("D31").Value = "1.175"
$ws.Range("E31").Value = "  -1.00%  "
$ws.Range("D32").Value = "4.078"
$ws.Range("E32").Value = "  +0.60%  "
$ws.Range("D33").Value = "5.797"
$ws.Range("E33").Value = "  +3.78%  "
$ws.Range("D34").Value = "0.08423"
$ws.Range("E34").Value = "  +0.97%  "
$ws.Range("D35").Value = "1.638"
$ws.Range("E35").Value = "  -2.56%  "
$ws.Range("D36").Value = "12.27"
$ws.Range("E36").Value = "  -0.60%  "
$ws.Range("D37").Value = "5.154"
$ws.Range("E37").Value = "  -0.90%  "
$ws.Range("D38").Value = "0.06035"
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "1.229"
$ws.Range("E39").Value = "  +2.10%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.02235"
$ws.Range("E40").Value = "  +2.03%  "
$ws.Range("E41").Value = "  +1.57%  "
$ws.Range("D42").Value = "8.227"
$ws.Range("E42").Value = "  +0.89%  "
$ws.Range("D43").Value = "0.9975"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").Value = "0.5945"
$ws.Range("E44").Value = "  +1.06%  "
$ws.Range("D45").Value = "13.66"
$ws.Range("E45").Value = "  +8.02%  "
$ws.Range("D46").Value = "3.858"
$ws.Range("E46").Value = "  +3.39%  "
$ws.Range("D47").Value = "0.5727"
$ws.Range("E47").Value = "  +2.45%  "
$ws.Range("D48").Value = "124.65"
$ws.Range("E48").Value = "  +2.11%  "
$ws.Range("D49").Value = "1.961"
$ws.Range("E49").Value = "  +1.34%  "
$ws.Range("D50").Value = "0.07008"
$ws.Range("E50").Value = "  +1.46%  "
$ws.Range("D51").Value = "1.192"
$ws.Range("E51").Value = "  +3.69%  "
